# example.calls.xlsx update:
#  - append a new data row (row 10) for patient Pt0 / sample Pt0_Bcells / chrY
#  - recompute VAF with the same formula pattern used by the other rows
#  - move the active selection to C12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Pt0"
$ws.Range("B10").Value = "Pt0_Bcells"
$ws.Range("C10").Value = "chrY"
$ws.Range("D10").Value = 1043
$ws.Range("E10").Value = 41
$ws.Range("F10").Value = 9
$ws.Range("G10").Formula = "=(F10/(F10+E10))"

$ws.Range("C12").Select() | Out-Null
